$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fase de Grupos")
$ws.Activate()

$ws.Range("F21").Value = 2
$ws.Range("H21").Value = 0

$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 1

$ws.Range("F33").Value = 0
$ws.Range("H33").Value = 0

$ws.Range("F34").Value = 0
$ws.Range("H34").Value = 1
